{"js": "// The \"\u041e\u0422\u0417\u042b\u0412\" (reviewer's feedback) cover sheet states the diploma-project\n// title. The title text is updated from\n//   \"...\u044d\u043b\u0435\u043a\u0442\u0440\u043e\u043d\u043d\u044b\u043c\u0438 \u0434\u043e\u043a\u0443\u043c\u0435\u043d\u0442\u0430\u043c\u0438 \u043d\u0430 \u043f\u0440\u043e\u0438\u0437\u0432\u043e\u0434\u0441\u0442\u0432\u0435\"  (at the production site)\n// to\n//   \"...\u044d\u043b\u0435\u043a\u0442\u0440\u043e\u043d\u043d\u044b\u043c\u0438 \u0434\u043e\u043a\u0443\u043c\u0435\u043d\u0442\u0430\u043c\u0438 \u043d\u0430 \u043f\u0440\u0435\u0434\u043f\u0440\u0438\u044f\u0442\u0438\u0438\"   (at the enterprise)\n// Do a targeted, case-sensitive search/replace for the exact word so only\n// the intended occurrence is touched.\nconst searchResults = context.document.body.search(\"\u043f\u0440\u043e\u0438\u0437\u0432\u043e\u0434\u0441\u0442\u0432\u0435\", {\n  matchCase: true,\n  matchWholeWord: true\n});\nsearchResults.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < searchResults.items.length; i++) {\n  searchResults.items[i].insertText(\"\u043f\u0440\u0435\u0434\u043f\u0440\u0438\u044f\u0442\u0438\u0438\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# The \"\u041e\u0422\u0417\u042b\u0412\" (reviewer's feedback) cover sheet states the diploma-project\n# title. The title text is updated from\n#   \"...\u044d\u043b\u0435\u043a\u0442\u0440\u043e\u043d\u043d\u044b\u043c\u0438 \u0434\u043e\u043a\u0443\u043c\u0435\u043d\u0442\u0430\u043c\u0438 \u043d\u0430 \u043f\u0440\u043e\u0438\u0437\u0432\u043e\u0434\u0441\u0442\u0432\u0435\"  (at the production site)\n# to\n#   \"...\u044d\u043b\u0435\u043a\u0442\u0440\u043e\u043d\u043d\u044b\u043c\u0438 \u0434\u043e\u043a\u0443\u043c\u0435\u043d\u0442\u0430\u043c\u0438 \u043d\u0430 \u043f\u0440\u0435\u0434\u043f\u0440\u0438\u044f\u0442\u0438\u0438\"   (at the enterprise)\n# Perform a targeted Find & Replace for the exact word so only the intended\n# occurrence in the title is touched.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"\u043f\u0440\u043e\u0438\u0437\u0432\u043e\u0434\u0441\u0442\u0432\u0435\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $true\n$find.MatchWildcards = $false\n\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"\u043f\u0440\u0435\u0434\u043f\u0440\u0438\u044f\u0442\u0438\u0438\"\n\n$find.Execute(\n    $find.Text,             # FindText\n    $true,                  # MatchCase\n    $true,                  # MatchWholeWord\n    $false,                 # MatchWildcards\n    $false,                 # MatchSoundsLike\n    $false,                 # MatchAllWordForms\n    $true,                  # Forward\n    [Microsoft.Office.Interop.Word.WdFindWrap]::wdFindContinue,\n    $false,                 # Format\n    $find.Replacement.Text, # ReplaceWith\n    [Microsoft.Office.Interop.Word.WdReplace]::wdReplaceAll\n)\n"}
